$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new values for column C (SearchRecord / ZOR)
$ws.Range("C1").Value = "SearchRecord"
$ws.Range("C2").Value = "ZOR"

# Update selection to match the diff (active cell C1)
$ws.Range("C1").Select()
